$d = $word.ActiveDocument

# --- Change 1: the law has been published already - update the opening wording ---
$old1 = "Dne 07.06.2023 Prezident ČR podepsal Novelu o Azilu, podle které"
$new1 = "Dne 21.06.2023 ve Sbírce zákonů v částce 87 pod číslem 173/2023 Sb. byla vyhlášena Novela z. o Azylu, podle které"

$range1 = $d.Content
$found1 = $range1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "Change1 found: $found1"

# --- Change 2: merge the run boundary between "Ministerstvo Vnitra " and "neomezovalo " ---
# (visible text is unchanged - only the underlying run split moves)
$old2 = "Ministerstvo Vnitra neomezovalo "
$range2 = $d.Content
$found2 = $range2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2)
Write-Output "Change2 found: $found2"

# --- Change 3: merge the run boundary between " " and "k " (around "povolení ... pobytu") ---
$old3 = " k "
$range3 = $d.Content
$found3 = $range3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)
Write-Output "Change3 found: $found3"
